$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (name/link rotation + refreshed price/volume figures)
$changes = @(
    @{ Cell = 'D2'; Value = '304.55'; Numeric = $true }
    @{ Cell = 'E2'; Value = '-0.23%'; Numeric = $true }
    @{ Cell = 'D3'; Value = '35.60'; Numeric = $true }
    @{ Cell = 'E3'; Value = '-3.92%'; Numeric = $true }
    @{ Cell = 'D4'; Value = '5.056'; Numeric = $true }
    @{ Cell = 'E4'; Value = '0.85%'; Numeric = $true }
    @{ Cell = 'D5'; Value = '0.07891'; Numeric = $true }
    @{ Cell = 'E5'; Value = '0.05%'; Numeric = $true }
    @{ Cell = 'D6'; Value = '2.119'; Numeric = $true }
    @{ Cell = 'E6'; Value = '-3.89%'; Numeric = $true }
    @{ Cell = 'B7'; Value = 'KuCoinToken'; Numeric = $false }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'; Numeric = $false }
    @{ Cell = 'D7'; Value = '7.906'; Numeric = $true }
    @{ Cell = 'E7'; Value = '-1.32%'; Numeric = $true }
    @{ Cell = 'B8'; Value = 'MXToken'; Numeric = $false }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Numeric = $false }
    @{ Cell = 'D8'; Value = '0.9233'; Numeric = $true }
    @{ Cell = 'E8'; Value = '0.03%'; Numeric = $true }
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange'; Numeric = $false }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; Numeric = $false }
    @{ Cell = 'D9'; Value = '0.09659'; Numeric = $true }
    @{ Cell = 'E9'; Value = '0.19%'; Numeric = $true }
    @{ Cell = 'B10'; Value = 'WazirX'; Numeric = $false }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; Numeric = $false }
    @{ Cell = 'D10'; Value = '0.1850'; Numeric = $true }
    @{ Cell = 'E10'; Value = '-2.18%'; Numeric = $true }
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken'; Numeric = $false }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; Numeric = $false }
    @{ Cell = 'D11'; Value = '0.08683'; Numeric = $true }
    @{ Cell = 'E11'; Value = '1.23%'; Numeric = $true }
    @{ Cell = 'B12'; Value = 'BitrueCoin'; Numeric = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; Numeric = $false }
    @{ Cell = 'D12'; Value = '0.03558'; Numeric = $true }
    @{ Cell = 'E12'; Value = '-3.52%'; Numeric = $true }
    @{ Cell = 'B13'; Value = 'BitMartToken'; Numeric = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; Numeric = $false }
    @{ Cell = 'D13'; Value = '0.09899'; Numeric = $true }
    @{ Cell = 'E13'; Value = '-0.78%'; Numeric = $true }
    @{ Cell = 'B14'; Value = 'BitForexToken'; Numeric = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; Numeric = $false }
    @{ Cell = 'D14'; Value = '0.001435'; Numeric = $true }
    @{ Cell = 'E14'; Value = '-2.96%'; Numeric = $true }
    @{ Cell = 'B15'; Value = 'TigerCash'; Numeric = $false }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'; Numeric = $false }
    @{ Cell = 'D15'; Value = '0.005669'; Numeric = $true }
    @{ Cell = 'E15'; Value = '0.55%'; Numeric = $true }
    @{ Cell = 'B16'; Value = 'LEO'; Numeric = $false }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; Numeric = $false }
    @{ Cell = 'D16'; Value = '3.472'; Numeric = $true }
    @{ Cell = 'E16'; Value = '0.65%'; Numeric = $true }
    @{ Cell = 'B17'; Value = 'GateToken'; Numeric = $false }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; Numeric = $false }
    @{ Cell = 'D17'; Value = '4.121'; Numeric = $true }
    @{ Cell = 'E17'; Value = '2.51%'; Numeric = $true }
    @{ Cell = 'E18'; Value = '17.38%'; Numeric = $true }
    @{ Cell = 'D19'; Value = '0.3370'; Numeric = $true }
    @{ Cell = 'E19'; Value = '-1.23%'; Numeric = $true }
    @{ Cell = 'D20'; Value = '0.1339'; Numeric = $true }
    @{ Cell = 'E20'; Value = '1.57%'; Numeric = $true }
    @{ Cell = 'D21'; Value = '5.160'; Numeric = $true }
    @{ Cell = 'E21'; Value = '8.48%'; Numeric = $true }
    @{ Cell = 'D22'; Value = '0.2210'; Numeric = $true }
    @{ Cell = 'E22'; Value = '0.45%'; Numeric = $true }
    @{ Cell = 'D23'; Value = '0.04504'; Numeric = $true }
    @{ Cell = 'E23'; Value = '-1.17%'; Numeric = $true }
    @{ Cell = 'E24'; Value = '-0.02%'; Numeric = $true }
    @{ Cell = 'D25'; Value = '0.004859'; Numeric = $true }
    @{ Cell = 'E25'; Value = '8.68%'; Numeric = $true }
    @{ Cell = 'E26'; Value = '-7.00%'; Numeric = $true }
    @{ Cell = 'D27'; Value = '0.0004755'; Numeric = $true }
    @{ Cell = 'E27'; Value = '0.10%'; Numeric = $true }
    @{ Cell = 'D39'; Value = '0.01838'; Numeric = $true }
    @{ Cell = 'E39'; Value = '-0.08%'; Numeric = $true }
    @{ Cell = 'D40'; Value = '0.04739'; Numeric = $true }
    @{ Cell = 'E40'; Value = '-0.64%'; Numeric = $true }
    @{ Cell = 'D41'; Value = '0.007877'; Numeric = $true }
    @{ Cell = 'E41'; Value = '-3.23%'; Numeric = $true }
    @{ Cell = 'D42'; Value = '0.1392'; Numeric = $true }
    @{ Cell = 'E42'; Value = '-0.50%'; Numeric = $true }
    @{ Cell = 'D43'; Value = '0.007758'; Numeric = $true }
    @{ Cell = 'E43'; Value = '2.58%'; Numeric = $true }
    @{ Cell = 'D44'; Value = '0.002192'; Numeric = $true }
    @{ Cell = 'E44'; Value = '-1.25%'; Numeric = $true }
    @{ Cell = 'D45'; Value = '0.01116'; Numeric = $true }
    @{ Cell = 'E45'; Value = '11.02%'; Numeric = $true }
    @{ Cell = 'D46'; Value = '0.00006301'; Numeric = $true }
    @{ Cell = 'E46'; Value = '0.31%'; Numeric = $true }
    @{ Cell = 'D47'; Value = '0.00000000751'; Numeric = $true }
    @{ Cell = 'E47'; Value = '0.13%'; Numeric = $true }
    @{ Cell = 'E48'; Value = '0.16%'; Numeric = $true }
    @{ Cell = 'E49'; Value = '76.86%'; Numeric = $true }
    @{ Cell = 'D50'; Value = '0.001903'; Numeric = $true }
    @{ Cell = 'E50'; Value = '10.54%'; Numeric = $true }
    @{ Cell = 'D51'; Value = '0.00002103'; Numeric = $true }
    @{ Cell = 'E51'; Value = '0.13%'; Numeric = $true }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    if ($change.Numeric) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $change.Value
}

